$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 5110
$ws.Range("J97").Value = 5110
$ws.Range("L97").Value = 15330
$ws.Range("N97").Value = -16322
$ws.Range("H116").Value = 6314.24
$ws.Range("I116").Value = 9370
$ws.Range("J116").Value = 3003.8333
$ws.Range("K116").Value = 9370
$ws.Range("L116").Value = 3003.8333
$ws.Range("M116").Value = -5928
$ws.Range("N116").Value = -9887.8333
$ws.Range("H132").Value = 766.0146999999999
$ws.Range("I132").Value = 692.069
$ws.Range("J132").Value = 1194.9
$ws.Range("K132").Value = 2076.207
$ws.Range("L132").Value = 3584.7
$ws.Range("M132").Value = 453.7930000000001
$ws.Range("N132").Value = -8644.700000000001
$ws.Range("H137").Value = 1929.1212
$ws.Range("I137").Value = 1589.8334
$ws.Range("J137").Value = 2833.889
$ws.Range("K137").Value = 4769.5002
$ws.Range("L137").Value = 8501.667000000001
$ws.Range("M137").Value = -2219.5002
$ws.Range("N137").Value = -13601.667
$ws.Range("H141").Value = 1376.4108
$ws.Range("I141").Value = 940.87805
$ws.Range("J141").Value = 2566.8667
$ws.Range("K141").Value = 2822.63415
$ws.Range("L141").Value = 7700.6001
$ws.Range("M141").Value = 2357.36585
$ws.Range("N141").Value = -18060.6001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4713.89
$ws.Range("I32").Value = 3461.2073
$ws.Range("J32").Value = 10420.556
$ws.Range("K32").Value = 3461.2073
$ws.Range("L32").Value = 10420.556
$ws.Range("M32").Value = -3174.2073
$ws.Range("N32").Value = -10994.556
$ws.Range("H61").Value = 3686.372
$ws.Range("I61").Value = 4037.5293
$ws.Range("J61").Value = 2359.7778
$ws.Range("K61").Value = 4037.5293
$ws.Range("L61").Value = 2359.7778
$ws.Range("M61").Value = -3825.5293
$ws.Range("N61").Value = -2783.7778
$ws.Range("H74").Value = 1154.4127
$ws.Range("I74").Value = 883.01886
$ws.Range("J74").Value = 2592.8
$ws.Range("K74").Value = 883.01886
$ws.Range("L74").Value = 2592.8
$ws.Range("M74").Value = -9.018860000000018
$ws.Range("N74").Value = -4340.8
$ws.Range("H77").Value = 1154.4127
$ws.Range("I77").Value = 883.01886
$ws.Range("J77").Value = 2592.8
$ws.Range("K77").Value = 4415.0943
$ws.Range("L77").Value = 12964
$ws.Range("M77").Value = -47.09429999999975
$ws.Range("N77").Value = -21700
$ws.Range("H122").Value = 801888.1
$ws.Range("I122").Value = 1068936.5
$ws.Range("J122").Value = 743
$ws.Range("K122").Value = 3206809.5
$ws.Range("L122").Value = 2229
$ws.Range("M122").Value = -3204359.5
$ws.Range("N122").Value = -7129
$ws.Range("H136").Value = 3686.372
$ws.Range("I136").Value = 4037.5293
$ws.Range("J136").Value = 2359.7778
$ws.Range("K136").Value = 12112.5879
$ws.Range("L136").Value = 7079.3334
$ws.Range("M136").Value = -9562.5879
$ws.Range("N136").Value = -12179.3334

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3803.1187
$ws.Range("I134").Value = 5119.857
$ws.Range("J134").Value = 2613.8064
$ws.Range("K134").Value = 15359.571
$ws.Range("L134").Value = 7841.4192
$ws.Range("M134").Value = -12824.571
$ws.Range("N134").Value = -12911.4192

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 223873.36
$ws.Range("I31").Value = 1346.3529
$ws.Range("J31").Value = 717302.8
$ws.Range("K31").Value = 1346.3529
$ws.Range("L31").Value = 717302.8
$ws.Range("M31").Value = -1051.3529
$ws.Range("N31").Value = -717892.8
$ws.Range("H34").Value = 223873.36
$ws.Range("I34").Value = 1346.3529
$ws.Range("J34").Value = 717302.8
$ws.Range("K34").Value = 1346.3529
$ws.Range("L34").Value = 717302.8
$ws.Range("M34").Value = -1144.3529
$ws.Range("N34").Value = -717706.8
$ws.Range("H58").Value = 1154.5938
$ws.Range("I58").Value = 606.4878
$ws.Range("J58").Value = 2131.652
$ws.Range("K58").Value = 606.4878
$ws.Range("L58").Value = 2131.652
$ws.Range("M58").Value = -403.4878
$ws.Range("N58").Value = -2537.652
$ws.Range("H122").Value = 1324095.1
$ws.Range("I122").Value = 2315980.2
$ws.Range("J122").Value = 1581.7778
$ws.Range("K122").Value = 6947940.600000001
$ws.Range("L122").Value = 4745.3334
$ws.Range("M122").Value = -6945490.600000001
$ws.Range("N122").Value = -9645.3334
$ws.Range("H132").Value = 1709.8718
$ws.Range("I132").Value = 1262.9508
$ws.Range("J132").Value = 3313.5293
$ws.Range("K132").Value = 3788.8524
$ws.Range("L132").Value = 9940.5879
$ws.Range("M132").Value = -1258.8524
$ws.Range("N132").Value = -15000.5879
$ws.Range("H134").Value = 1628.2113
$ws.Range("I134").Value = 1849.326
$ws.Range("J134").Value = 1221.36
$ws.Range("K134").Value = 5547.978
$ws.Range("L134").Value = 3664.08
$ws.Range("M134").Value = -3012.978
$ws.Range("N134").Value = -8734.08
$ws.Range("H136").Value = 1154.5938
$ws.Range("I136").Value = 606.4878
$ws.Range("J136").Value = 2131.652
$ws.Range("K136").Value = 1819.4634
$ws.Range("L136").Value = 6394.956
$ws.Range("M136").Value = 730.5365999999999
$ws.Range("N136").Value = -11494.956

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1486.5555
$ws.Range("J34").Value = 2174.75
$ws.Range("L34").Value = 6524.25
$ws.Range("N34").Value = -6692.25
$ws.Range("H39").Value = 3333.3333
$ws.Range("J39").Value = 3333.3333
$ws.Range("L39").Value = 9999.999899999999
$ws.Range("N39").Value = -10587.9999
$ws.Range("H55").Value = 1809.7273
$ws.Range("J55").Value = 1809.7273
$ws.Range("L55").Value = 5429.1819
$ws.Range("N55").Value = -5783.1819
$ws.Range("H92").Value = 803
$ws.Range("J92").Value = 803
$ws.Range("L92").Value = 2409
$ws.Range("N92").Value = -4905
$ws.Range("H98").Value = 20000420
$ws.Range("I98").Value = 250
$ws.Range("J98").Value = 33333866
$ws.Range("K98").Value = 750
$ws.Range("L98").Value = 100001598
$ws.Range("M98").Value = 748
$ws.Range("N98").Value = -100004594
$ws.Range("H107").Value = 334.90244
$ws.Range("I107").Value = 215.375
$ws.Range("J107").Value = 363.87878
$ws.Range("K107").Value = 646.125
$ws.Range("L107").Value = 1091.63634
$ws.Range("M107").Value = 1273.875
$ws.Range("N107").Value = -4931.63634
$ws.Range("H122").Value = 2517.681
$ws.Range("I122").Value = 376.94736
$ws.Range("J122").Value = 3970.3215
$ws.Range("K122").Value = 3392.52624
$ws.Range("L122").Value = 35732.8935
$ws.Range("M122").Value = -942.5262400000001
$ws.Range("N122").Value = -40632.8935
$ws.Range("H131").Value = 2041623.9
$ws.Range("I131").Value = 7143136
$ws.Range("J131").Value = 1018.88574
$ws.Range("K131").Value = 21429408
$ws.Range("L131").Value = 3056.65722
$ws.Range("M131").Value = -21424368
$ws.Range("N131").Value = -13136.65722
$ws.Range("H132").Value = 10058232
$ws.Range("I132").Value = 2875.8
$ws.Range("J132").Value = 16342829
$ws.Range("K132").Value = 25882.2
$ws.Range("L132").Value = 147085461
$ws.Range("M132").Value = -23352.2
$ws.Range("N132").Value = -147090521

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 5800
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("H132").Value = 2304.4
$ws.Range("I132").Value = 1853.9354
$ws.Range("J132").Value = 3301.8572
$ws.Range("K132").Value = 5561.8062
$ws.Range("L132").Value = 9905.571599999999
$ws.Range("M132").Value = -3031.8062
$ws.Range("N132").Value = -14965.5716

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 30068.857
$ws.Range("I7").Value = 44947.176
$ws.Range("J7").Value = 1552.0834
$ws.Range("K7").Value = 44947.176
$ws.Range("L7").Value = 1552.0834
$ws.Range("M7").Value = -44835.176
$ws.Range("N7").Value = -1776.0834
$ws.Range("H40").Value = 30305232
$ws.Range("I40").Value = 38463070
$ws.Range("J40").Value = 4683.4287
$ws.Range("K40").Value = 38463070
$ws.Range("L40").Value = 4683.4287
$ws.Range("M40").Value = -38462934
$ws.Range("N40").Value = -4955.4287
$ws.Range("H82").Value = 2206896
$ws.Range("I82").Value = 5001350
$ws.Range("J82").Value = 343926.66
$ws.Range("K82").Value = 5001350
$ws.Range("L82").Value = 343926.66
$ws.Range("M82").Value = -5000989
$ws.Range("N82").Value = -344648.66
$ws.Range("H85").Value = 2206896
$ws.Range("I85").Value = 5001350
$ws.Range("J85").Value = 343926.66
$ws.Range("K85").Value = 5001350
$ws.Range("L85").Value = 343926.66
$ws.Range("M85").Value = -5000102
$ws.Range("N85").Value = -346422.66
$ws.Range("H126").Value = 30068.857
$ws.Range("I126").Value = 44947.176
$ws.Range("J126").Value = 1552.0834
$ws.Range("K126").Value = 134841.528
$ws.Range("L126").Value = 4656.2502
$ws.Range("M126").Value = -132371.528
$ws.Range("N126").Value = -9596.2502
$ws.Range("H132").Value = 14065936
$ws.Range("I132").Value = 17241150
$ws.Range("J132").Value = 4271.143
$ws.Range("K132").Value = 51723450
$ws.Range("L132").Value = 12813.429
$ws.Range("M132").Value = -51720920
$ws.Range("N132").Value = -17873.429
$ws.Range("H136").Value = 8402.121999999999
$ws.Range("I136").Value = 5421.7095
$ws.Range("J136").Value = 17641.4
$ws.Range("K136").Value = 16265.1285
$ws.Range("L136").Value = 52924.2
$ws.Range("M136").Value = -13715.1285
$ws.Range("N136").Value = -58024.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3676.1428
$ws.Range("I122").Value = 3622.1667
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 10866.5001
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -8416.500100000001
$ws.Range("N122").Value = -16900
$ws.Range("H132").Value = 12619.667
$ws.Range("I132").Value = 14628.31
$ws.Range("J132").Value = 1649.3846
$ws.Range("K132").Value = 43884.93
$ws.Range("L132").Value = 4948.1538
$ws.Range("M132").Value = -41354.93
$ws.Range("N132").Value = -10008.1538
$ws.Range("H136").Value = 9262146
$ws.Range("J136").Value = 20002020
$ws.Range("L136").Value = 60006060
$ws.Range("N136").Value = -60011160

# Clear M33 on GSM (cell removed entirely per diff)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M33").ClearContents()

Write-Output "Applied all Ixion Profits updates"